$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text (affects every cell using this shared string: C2, C3, C4)
$oldText = "PC DELL ALL IN ONE CORE I7"
$newText = "PC DELL ALL IN ONE CORE I9 mem 32gb"

$usedRange = $ws.UsedRange
$foundCell = $usedRange.Find($oldText)
if ($foundCell -ne $null) {
    $firstAddress = $foundCell.Address()
    do {
        $foundCell.Value = $newText
        $foundCell = $usedRange.FindNext($foundCell)
    } while ($foundCell -ne $null -and $foundCell.Address() -ne $firstAddress)
}

# Update the active selection to C3:C4 with C3 as the active cell
$ws.Range("C4:C3").Select()
